$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuickHour")

# --- Row 3 header cells: rename BT/PD routes with fuller descriptions ---
$ws.Range("T3").Value = "BT115" + [char]10 + "PARIS DAUPHINE"
$ws.Range("U3").Value = "BT131 REPRISE"
$ws.Range("V3").Value = "BT134" + [char]10 + "LAYANI "
$ws.Range("W3").Value = "BT144 J.D'ARC"
$ws.Range("X3").Value = "BT159 "
$ws.Range("Y3").Value = "BT160" + [char]10 + "UFR " + [char]10 + "ANGERS"
$ws.Range("Z3").Value = "BT163 TOULOUSE"
$ws.Range("AA3").Value = "BT171 BESSIER"
$ws.Range("AB3").Value = "BT172" + [char]10 + "TOUR NEPTUNE"
$ws.Range("AH3").Value = "PD102 LAYANI"
$ws.Range("AI3").Value = "PD106 TRAPPES"
$ws.Range("AJ3").Value = "PD117" + [char]10 + "CARAC "
$ws.Range("AK3").Value = "PD13 NANTES"
$ws.Range("AL3").Value = "PD16" + [char]10 + "LA RAPÉE"
$ws.Range("AM3").Value = "PD18 GUYANCOURT"
$ws.Range("AN3").Value = "PD36" + [char]10 + "R. ALBERT"
$ws.Range("AO3").Value = "PD38 BESSIER"
$ws.Range("AP3").Value = "PD39 LIMOURS"
$ws.Range("AQ3").Value = "PD40 BESSIER"
$ws.Range("AR3").Value = "PD41 BESSIER"
$ws.Range("AS3").Value = "PD42 BESSIER"
$ws.Range("AT3").Value = "PD466 COMMYNES"
$ws.Range("AU3").Value = "PD467 COMMYNES"
$ws.Range("AV3").Value = "PD58 CLÉRET"
$ws.Range("AX3").Value = "X"
$ws.Range("AY3").Value = "X*"
$ws.Range("AZ3").Value = "X**"
$ws.Range("BA3").Value = "X***"
$ws.Range("BB3").Value = "X****"
$ws.Range("BF3").Value = "BT115"
$ws.Range("BJ3").Value = "BT131"
$ws.Range("BK3").Value = "BT134"
$ws.Range("BM3").Value = "BT144"
$ws.Range("BW3").Value = "BT159"
$ws.Range("BX3").Value = "BT160"
$ws.Range("BZ3").Value = "BT163"
$ws.Range("CD3").Value = "BT171"
$ws.Range("CE3").Value = "BT172"

# --- Row 4 totals: the newly-used X/X*/X**/X***/X**** columns now carry a 0 total like their neighbours ---
$ws.Range("AX4").Value = 0
$ws.Range("AY4").Value = 0
$ws.Range("AZ4").Value = 0
$ws.Range("BA4").Value = 0
$ws.Range("BB4").Value = 0
